# Sheets.xlsx scheduled refresh: update cached market-price / profit columns (H-N)
# for the affected leve rows across all 8 job sheets. Values below come from
# a Universalis price-data refresh; no formulas are involved (cells are static).

$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 15: Morning Glass of Ether / Ether
$ws.Range("H15").Value = 626236.5600000001
$ws.Range("I15").Value = 626236.5600000001
$ws.Range("K15").Value = 1878709.68
$ws.Range("M15").Value = -1878540.68

# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 634
$ws.Range("I28").Value = 556.4545000000001
$ws.Range("K28").Value = 556.4545000000001
$ws.Range("M28").Value = -71.45450000000005

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 9673.308000000001
$ws.Range("J62").Value = 9861.666999999999
$ws.Range("L62").Value = 9861.666999999999
$ws.Range("N62").Value = -11109.667

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 9673.308000000001
$ws.Range("J65").Value = 9861.666999999999
$ws.Range("L65").Value = 49308.335
$ws.Range("N65").Value = -55548.335

# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 10348.223
$ws.Range("J69").Value = 11934.5
$ws.Range("L69").Value = 35803.5
$ws.Range("N69").Value = -37551.5

# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 10348.223
$ws.Range("J72").Value = 11934.5
$ws.Range("L72").Value = 107410.5
$ws.Range("N72").Value = -116146.5

# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 1166.3478
$ws.Range("I92").Value = 1128.6
$ws.Range("J92").Value = 1418
$ws.Range("K92").Value = 1128.6
$ws.Range("L92").Value = 1418
$ws.Range("M92").Value = 119.4000000000001
$ws.Range("N92").Value = -3914

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 36289.47
$ws.Range("I98").Value = 38595
$ws.Range("J98").Value = 18998
$ws.Range("K98").Value = 38595
$ws.Range("L98").Value = 18998
$ws.Range("M98").Value = -37097
$ws.Range("N98").Value = -21994

# Row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 1439.2632
$ws.Range("I111").Value = 1387
$ws.Range("J111").Value = 1585.6
$ws.Range("K111").Value = 4161
$ws.Range("L111").Value = 4756.799999999999
$ws.Range("M111").Value = -1094
$ws.Range("N111").Value = -10890.8

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 36289.47
$ws.Range("I122").Value = 38595
$ws.Range("J122").Value = 18998
$ws.Range("K122").Value = 115785
$ws.Range("L122").Value = 56994
$ws.Range("M122").Value = -113335
$ws.Range("N122").Value = -61894

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 882343.75
$ws.Range("I137").Value = 1161633.8
$ws.Range("K137").Value = 3484901.4
$ws.Range("M137").Value = -3482351.4

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4301.56
$ws.Range("I138").Value = 1460.48
$ws.Range("J138").Value = 5248.5864
$ws.Range("K138").Value = 4381.440000000001
$ws.Range("L138").Value = 15745.7592
$ws.Range("M138").Value = 758.5599999999995
$ws.Range("N138").Value = -26025.7592

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 9139.4
$ws.Range("I141").Value = 9935.637000000001
$ws.Range("J141").Value = 6949.75
$ws.Range("K141").Value = 29806.911
$ws.Range("L141").Value = 20849.25
$ws.Range("M141").Value = -24626.911
$ws.Range("N141").Value = -31209.25


# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 8620.552
$ws.Range("I32").Value = 8337.508
$ws.Range("K32").Value = 8337.508
$ws.Range("M32").Value = -8050.508

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1703.0857
$ws.Range("I74").Value = 764.2727
$ws.Range("K74").Value = 764.2727
$ws.Range("M74").Value = 109.7273

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1703.0857
$ws.Range("I77").Value = 764.2727
$ws.Range("K77").Value = 3821.3635
$ws.Range("M77").Value = 546.6365000000001

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 2189.3125
$ws.Range("I110").Value = 1855.2273
$ws.Range("K110").Value = 1855.2273
$ws.Range("M110").Value = 189.7727

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 603502.75
$ws.Range("I122").Value = 3292.7297
$ws.Range("K122").Value = 9878.1891
$ws.Range("M122").Value = -7428.1891

# Row 129: In-kweh-dible Cooking / Manganese Chocobo Frypan
$ws.Range("H129").Value = 59561.25
$ws.Range("J129").Value = 59500
$ws.Range("L129").Value = 59500
$ws.Range("N129").Value = -69500


# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 4030.0952
$ws.Range("I3").Value = 3860.7646
$ws.Range("J3").Value = 4749.75
$ws.Range("K3").Value = 3860.7646
$ws.Range("L3").Value = 4749.75
$ws.Range("M3").Value = -3746.7646
$ws.Range("N3").Value = -4977.75

# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 9384.6
$ws.Range("I20").Value = 7669.4
$ws.Range("K20").Value = 7669.4
$ws.Range("M20").Value = -7422.4


# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 907
$ws.Range("I22").Value = 512.7143
$ws.Range("K22").Value = 512.7143
$ws.Range("M22").Value = -162.7143

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2932.7778
$ws.Range("I31").Value = 1852.1052
$ws.Range("K31").Value = 1852.1052
$ws.Range("M31").Value = -1557.1052

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2932.7778
$ws.Range("I34").Value = 1852.1052
$ws.Range("K34").Value = 1852.1052
$ws.Range("M34").Value = -1650.1052

# Row 59: Bow Down to Magic / Crab Bow
$ws.Range("H59").Value = 99127
$ws.Range("J59").Value = 99127
$ws.Range("L59").Value = 99127
$ws.Range("N59").Value = -101417

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 4811.3076
$ws.Range("I132").Value = 5461.136
$ws.Range("J132").Value = 1237.25
$ws.Range("K132").Value = 16383.408
$ws.Range("L132").Value = 3711.75
$ws.Range("M132").Value = -13853.408
$ws.Range("N132").Value = -8771.75

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1813.5758
$ws.Range("I134").Value = 1270.2069
$ws.Range("K134").Value = 3810.620699999999
$ws.Range("M134").Value = -1275.620699999999

# Row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 202155.66
$ws.Range("J141").Value = 216496.17
$ws.Range("L141").Value = 216496.17
$ws.Range("N141").Value = -226856.17


# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 56: Culture Club / Crowned Pie
$ws.Range("H56").Value = 7849.129
$ws.Range("I56").Value = 7849.129
$ws.Range("K56").Value = 7849.129
$ws.Range("M56").Value = -7319.129

# Row 81: It Goes Down Smoothly / Frozen Spirits
$ws.Range("H81").Value = 90039.25
$ws.Range("I81").Value = 334837.66
$ws.Range("J81").Value = 8439.777
$ws.Range("K81").Value = 1004512.98
$ws.Range("L81").Value = 25319.331
$ws.Range("M81").Value = -1003389.98
$ws.Range("N81").Value = -27565.331

# Row 84: Quenching the Flame (L) / Frozen Spirits
$ws.Range("H84").Value = 90039.25
$ws.Range("I84").Value = 334837.66
$ws.Range("J84").Value = 8439.777
$ws.Range("K84").Value = 3013538.94
$ws.Range("L84").Value = 75957.993
$ws.Range("M84").Value = -3007922.94
$ws.Range("N84").Value = -87189.993

# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 9469.799999999999
$ws.Range("I134").Value = 19535
$ws.Range("K134").Value = 58605
$ws.Range("M134").Value = -53535


# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 4: Arms for the Poor / Bone Brand
$ws.Range("H4").Value = 6933.3335
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 8531.344999999999
$ws.Range("I70").Value = 7047.4287
$ws.Range("K70").Value = 7047.4287
$ws.Range("M70").Value = -6777.4287

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 8531.344999999999
$ws.Range("I73").Value = 7047.4287
$ws.Range("K73").Value = 7047.4287
$ws.Range("M73").Value = -6111.4287

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 553.1053000000001
$ws.Range("I107").Value = 585.0714
$ws.Range("J107").Value = 463.6
$ws.Range("K107").Value = 585.0714
$ws.Range("L107").Value = 463.6
$ws.Range("M107").Value = 1334.9286
$ws.Range("N107").Value = -4303.6

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 11690.064
$ws.Range("I122").Value = 9834.643
$ws.Range("J122").Value = 13218.059
$ws.Range("K122").Value = 29503.929
$ws.Range("L122").Value = 39654.177
$ws.Range("M122").Value = -27053.929
$ws.Range("N122").Value = -44554.177

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2937.0256
$ws.Range("I132").Value = 2942.2727
$ws.Range("K132").Value = 8826.8181
$ws.Range("M132").Value = -6296.8181

# Row 134: Guaranteed Gem / Ihuykanite
$ws.Range("H134").Value = 80081
$ws.Range("J134").Value = 80081
$ws.Range("L134").Value = 240243
$ws.Range("N134").Value = -245313


# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 5761.037
$ws.Range("I16").Value = 6343.1816
$ws.Range("J16").Value = 3199.6
$ws.Range("K16").Value = 6343.1816
$ws.Range("L16").Value = 3199.6
$ws.Range("M16").Value = -6173.1816
$ws.Range("N16").Value = -3539.6

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 873.4783
$ws.Range("J22").Value = 973.6
$ws.Range("L22").Value = 973.6
$ws.Range("N22").Value = -1563.6

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 873.4783
$ws.Range("J27").Value = 973.6
$ws.Range("L27").Value = 973.6
$ws.Range("N27").Value = -1187.6

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 6892.364
$ws.Range("I122").Value = 7493.7
$ws.Range("K122").Value = 22481.1
$ws.Range("M122").Value = -20031.1

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 1360668.4
$ws.Range("I132").Value = 1496385.2
$ws.Range("K132").Value = 4489155.6
$ws.Range("M132").Value = -4486625.6

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 7154.1333
$ws.Range("I136").Value = 7332.6665
$ws.Range("K136").Value = 21997.9995
$ws.Range("M136").Value = -19447.9995


# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 20478.438
$ws.Range("I107").Value = 2629.2856
$ws.Range("K107").Value = 7887.8568
$ws.Range("M107").Value = -5967.8568

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1490.7646
$ws.Range("I113").Value = 887.3103599999999
$ws.Range("K113").Value = 2661.93108
$ws.Range("M113").Value = -491.9310799999998

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 4685.5483
$ws.Range("I122").Value = 2359.75
$ws.Range("K122").Value = 7079.25
$ws.Range("M122").Value = -4629.25

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 5289.0615
$ws.Range("I132").Value = 6218.82
$ws.Range("J132").Value = 2453.3
$ws.Range("K132").Value = 18656.46
$ws.Range("L132").Value = 7359.900000000001
$ws.Range("M132").Value = -16126.46
$ws.Range("N132").Value = -12419.9

